# Refresh the cryptos snapshot: updated Price (D) / Volume(1h) (E) figures
# for every listed coin, plus two rows (36/37 and 49/50) where the ranking
# order of two coins swapped, so the Coin name/Link/Price/Volume for those
# four rows are fully replaced.
#
# Price values are text in this sheet (e.g. "3.228.90" uses "." as a
# thousands separator, which is not a valid Excel number). Assigning a
# leading single-quote forces Excel to keep numeric-looking prices (like
# "0.999") as text instead of silently converting them to a number; the
# Style reset afterwards clears the resulting quote-prefix cell style so
# the cell's formatting is left exactly as it was before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.409.40'
$ws.Range('E2').Value = '  -7.01%  '

# Row 3
$ws.Range('D3').Value = '3.249.92'
$ws.Range('E3').Value = '  -9.50%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '''176.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -15.72%  '

# Row 6
$ws.Range('D6').Value = '''511.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -10.48%  '

# Row 7
$ws.Range('D7').Value = '''0.585'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.33%  '

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').Value = '3.242.41'
$ws.Range('E9').Value = '  -9.49%  '

# Row 10
$ws.Range('D10').Value = '''0.608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -10.63%  '

# Row 11
$ws.Range('D11').Value = '''56.88'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.51%  '

# Row 12
$ws.Range('D12').Value = '''0.129'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.88%  '

# Row 13
$ws.Range('D13').Value = '''0.0000251'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -10.74%  '

# Row 14
$ws.Range('D14').Value = '''9.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.93%  '

# Row 15
$ws.Range('D15').Value = '3.750.30'
$ws.Range('E15').Value = '  -9.62%  '

# Row 16
$ws.Range('D16').Value = '''0.118'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.69%  '

# Row 17
$ws.Range('D17').Value = '3.229.48'
$ws.Range('E17').Value = '  -9.70%  '

# Row 18
$ws.Range('D18').Value = '63.075.83'
$ws.Range('E18').Value = '  -7.17%  '

# Row 19
$ws.Range('D19').Value = '''16.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.64%  '

# Row 20
$ws.Range('D20').Value = '''10.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -12.54%  '

# Row 21
$ws.Range('D21').Value = '''0.933'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -12.45%  '

# Row 22
$ws.Range('D22').Value = '''366.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.24%  '

# Row 23
$ws.Range('D23').Value = '''78.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.97%  '

# Row 24
$ws.Range('D24').Value = '''10.89'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -12.20%  '

# Row 25
$ws.Range('D25').Value = '''3.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -13.84%  '

# Row 26
$ws.Range('D26').Value = '''5.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.13%  '

# Row 27
$ws.Range('D27').Value = '''3.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.12%  '

# Row 28
$ws.Range('D28').Value = '''2.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.53%  '

# Row 29
$ws.Range('D29').Value = '''11.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.26%  '

# Row 30
$ws.Range('D30').Value = '''8.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.93%  '

# Row 31
$ws.Range('D31').Value = '''640.47'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.22%  '

# Row 32
$ws.Range('D32').Value = '''28.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.24%  '

# Row 33
$ws.Range('D33').Value = '''6.59'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -14.03%  '

# Row 34
$ws.Range('D34').Value = '''11.01'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.40%  '

# Row 35
$ws.Range('D35').Value = '''58.67'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.14%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.102'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.31%  '

# Row 37
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '''1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.12%  '

# Row 38
$ws.Range('D38').Value = '''35.47'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -13.59%  '

# Row 39
$ws.Range('D39').Value = '''0.374'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.23%  '

# Row 40
$ws.Range('D40').Value = '''0.995'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.24%  '

# Row 41
$ws.Range('D41').Value = '''0.122'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.41%  '

# Row 42
$ws.Range('D42').Value = '2.845.28'
$ws.Range('E42').Value = '  -10.16%  '

# Row 43
$ws.Range('D43').Value = '0.0₃0647'
$ws.Range('E43').Value = '  -14.20%  '

# Row 44
$ws.Range('D44').Value = '''2.62'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -18.83%  '

# Row 45
$ws.Range('D45').Value = '''2.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.31%  '

# Row 46
$ws.Range('D46').Value = '''2.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -12.37%  '

# Row 47
$ws.Range('D47').Value = '''2.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.44%  '

# Row 48
$ws.Range('D48').Value = '''0.0377'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.71%  '

# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.122'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.20%  '

# Row 50
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '''2.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.49%  '

# Row 51
$ws.Range('D51').Value = '''131.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.24%  '

